$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 31 into the new row 32 so A32/B32/C32/D32 inherit the exact
# same values + types as row 31 currently has (in particular B32 stays a
# text "3", matching the pre-fix shape of B31).
$ws.Range("A31:H31").Copy($ws.Range("A32:H32"))

# Row 31: B31 was stored as text "3"; fix it to be a real number 3.
$ws.Range("B31").Value = 3

# Row 32 gets its own issue_type/id/source_file/text values.
$ws.Range("E32").Value = "OTH"
$ws.Range("F32").Value = "aa721c36-81b2-451c-915e-fe15286fe992"
$ws.Range("G32").Value = "SygwwGbRW_annotated.xlsx"
$ws.Range("H32").Value = "This is NOT a proper navigation agent."
